$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export data has moved forward by one day: the oldest row
# (2025-11-03) drops off, every remaining row shifts up by one, and two
# fresh rows (2026-01-31, 2026-02-01) land at the bottom.
#
# We shift via Copy/PasteSpecial (not `.Value = ...`) because the
# date column stores its dates as plain text (shared strings), and a
# direct `.Value = "2026-01-31"` assignment gets auto-coerced into a
# real date serial number by Excel's type inference. Copy/PasteSpecial
# duplicates the already-stored text cell verbatim, sidestepping that
# inference entirely.

# Shift rows 3..90 up into rows 2..89 (process top-down so each source
# row is read before it becomes a paste target).
for ($r = 2; $r -le 89; $r++) {
    $srcRow = $r + 1
    $ws.Range("A" + $srcRow + ":C" + $srcRow).Copy()
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial()
}

# Build the two new trailing date strings off-sheet (as text formula
# results) and paste-special just the values in, so they land as plain
# shared-string text instead of being reinterpreted as dates.
$scratch = $ws.Cells.Item(1, 10)

$scratch.Formula = "=""2026-01-31"""
$scratch.Copy()
$ws.Cells.Item(90, 1).PasteSpecial(-4163)

$scratch.Formula = "=""2026-02-01"""
$scratch.Copy()
$ws.Cells.Item(91, 1).PasteSpecial(-4163)

$scratch.ClearContents()

$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(90, 3).Value = 28

$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(91, 3).Value = 28
